# Commit message: "add the NA's under duplicate_image_filename"
#
# Column E ("duplicate_image_filename", header already present in E1) has
# no data for the trial rows below it. Fill rows 2-21 (the practice rows
# 2-5 and the main trial rows 6-21) with "NA".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E21").Value = "NA"
